$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Add row 17 to Logs sheet (new test mail entry) ---
$logs.Range("A17").Value = "Kunt u mij uitschrijven voor de nieuwsbrief?"
$logs.Range("B17").Value = "mailmind.test@zohomail.eu"
$logs.Range("C17").Value = "Testmail #17: Kunt u mij uitschrijven voor de nieuwsbrief?"
$logs.Range("D17").Value = "Afmelding / Nieuwsbrief"
$logs.Range("E17").Value = "Beste klant,`nDank voor uw bericht. Om u uit te schrijven voor onze nieuwsbrief, hebben wij uw e-mailadres nodig. Kunt u ons alstublieft het e-mailadres sturen waar u voor uitgeschreven wilt worden?`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F17").Value = "2025-07-22 12:48:20"
$logs.Range("G17").Value = "Ja"
$logs.Range("H17").Value = "Nee"
$logs.Range("I17").Value = "Ja"
$logs.Range("J17").Value = "Ja"

# The multi-line text in E17 makes the engine auto-size the new row's
# height; reset it back to the sheet's standard (non-custom) height so the
# row matches the other, untouched rows.
$logs.Rows.Item(17).AutoFit() | Out-Null
$logs.Rows.Item(17).UseStandardHeight = $true

# --- Extend conditional formatting ranges from row 16 to row 17 ---
$logs.Range("D2:D16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D17"))
$logs.Range("G2:G16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G17"))
$logs.Range("H2:H16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H17"))
$logs.Range("I2:I16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I17"))
$logs.Range("J2:J16").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J17"))

# --- Add row 9 to Dashboard sheet ---
$dash.Range("A9").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B9").Value = 1

# --- Update chart source ranges to include the new Dashboard row ---
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$9"
$series.Values = "='Dashboard'!`$B`$2:`$B`$9"
